$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.286.29"
$ws.Range("E2").Value = "  +1.95%  "

$ws.Range("D3").Value = "3.386.01"
$ws.Range("E3").Value = "  +1.67%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  +0.86%  "

$ws.Range("E9").Value = "  +6.08%  "

$ws.Range("E10").Value = "  +1.38%  "

$ws.Range("E11").Value = "  +3.03%  "

$ws.Range("E12").Value = "  +2.97%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "676.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.08%  "

$ws.Range("D15").Value = "3.932.39"
$ws.Range("E15").Value = "  +1.48%  "

$ws.Range("D16").Value = "69.292.20"
$ws.Range("E16").Value = "  +1.81%  "

$ws.Range("E17").Value = "  +1.84%  "

$ws.Range("D18").Value = "3.389.58"
$ws.Range("E18").Value = "  +1.78%  "

$ws.Range("E19").Value = "  +1.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.61%  "

$ws.Range("E21").Value = "  +0.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.16%  "

$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("E26").Value = "  +0.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.98%  "

$ws.Range("E29").Value = "  +1.38%  "

$ws.Range("E30").Value = "  -1.49%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "553.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.01%  "

$ws.Range("E34").Value = "  +0.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.99%  "

$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("D37").Value = "3.682.89"
$ws.Range("E37").Value = "  -0.64%  "

$ws.Range("E38").Value = "  +5.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.41%  "

$ws.Range("E41").Value = "  +0.84%  "

$ws.Range("D42").Value = "0.0₃0696"
$ws.Range("E42").Value = "  +3.12%  "

$ws.Range("E43").Value = "  +0.27%  "

$ws.Range("E44").Value = "  +3.28%  "

$ws.Range("E45").Value = "  -1.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("E47").Value = "  +0.55%  "

$ws.Range("E48").Value = "  +4.97%  "

$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.66%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.25%  "
